$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "extr1..extr8" block (rows 8-15) down by two rows,
# to rows 10-17, to make room for two new rows (line7, line8).
# Walk bottom-up so we never overwrite a row before reading it.
for ($r = 15; $r -ge 8; $r--) {
    $dest = $r + 2
    $ws.Cells.Item($dest, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($dest, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($dest, 3).Value = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($dest, 4).Value = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($dest, 5).Value = $ws.Cells.Item($r, 5).Value2
}

# Rows 16-17 are brand new territory on the sheet: give column A the same
# (bold/bordered/centered) formatting the rest of the index column uses.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)  # xlPasteFormats

# New row 8: line7
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# New row 9: line8
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# Renumber the "index" column (A) for the rows that shifted down
# (previously rows 8-15, now rows 10-17) to continue the sequence.
for ($r = 10; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
